# Update the "F" column (想去人数 / "number wanting to go") values on the
# "展览" and "全部类型" worksheets to reflect refreshed figures.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F, identical on both sheets.
$updates = @{
    2  = 325
    4  = 10416
    6  = 946
    7  = 24
    8  = 1286
    9  = 7202
    11 = 445
    12 = 206
    13 = 129
    14 = 3209
    16 = 316
    17 = 675
    22 = 1650
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
